# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" worksheet (fund-holding detail) right before
#   the "总计" (totals) sheet, formatted like the existing quarter sheets.
# - Add a corresponding summary row at the top of "总计" (shifting the
#   previously-existing rows down and renumbering the index column).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Clone the look (styles/borders/fonts) of the prior quarter sheet so the
# new tab matches the existing formatting, then overwrite with real data.
$wb.Worksheets.Item("2021-Q4").Range("A1:H13").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
  @(0,  "000697", "汇添富移动互联股票",               "24.17", "83.18", "4.24", "1.0248", 5),
  @(1,  "020003", "国泰金龙行业混合",                 "12.03", "69.69", "4.11", "0.4944", 8),
  @(2,  "160212", "国泰估值优势混合 (LOF)",            "8.98",  "62.69", "4.23", "0.3799", 5),
  @(3,  "013123", "汇添富精选核心优势一年持有混合A",    "6.15",  "66.61", "3.55", "0.2183", 7),
  @(4,  "460009", "华泰柏瑞量化先行混合A",             "9.13",  "90.47", "0.85", "0.0776", 9),
  @(5,  "006502", "财通集成电路产业股票A",             "1.29",  "79.76", "5.28", "0.0681", 4),
  @(6,  "006503", "财通集成电路产业股票C",             "0.46",  "79.76", "5.28", "0.0243", 4),
  @(7,  "002292", "诺安益鑫灵活配置混合",              "0.30",  "50.08", "4.25", "0.0128", 5),
  @(8,  "013124", "汇添富精选核心优势一年持有混合C",    "0.30",  "66.61", "3.55", "0.0106", 7),
  @(9,  "010246", "华泰柏瑞量化先行混合C",             "0.12",  "90.47", "0.85", "0.0010", 9),
  @(10, "004833", "先锋聚利灵活配置混合A",             "0.02",  "94.68", "4.92", "0.0010", 4),
  @(11, "004834", "先锋聚利灵活配置混合C",             "0.02",  "94.68", "4.92", "0.0010", 4)
)

$row = 2
foreach ($item in $fundRows) {
    $newSheet.Range("A$row").Value = $item[0]

    # Fund code / scale / position ratios are stored as plain text in the
    # source data (e.g. leading zeros in codes), so force text format
    # before assigning -- otherwise Excel coerces numeric-looking strings.
    $newSheet.Range("B$row").NumberFormat = "@"
    $newSheet.Range("B$row").Value = $item[1]

    $newSheet.Range("C$row").Value = $item[2]

    $newSheet.Range("D$row").NumberFormat = "@"
    $newSheet.Range("D$row").Value = $item[3]
    $newSheet.Range("E$row").NumberFormat = "@"
    $newSheet.Range("E$row").Value = $item[4]
    $newSheet.Range("F$row").NumberFormat = "@"
    $newSheet.Range("F$row").Value = $item[5]
    $newSheet.Range("G$row").NumberFormat = "@"
    $newSheet.Range("G$row").Value = $item[6]

    $newSheet.Range("H$row").Value = $item[7]

    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

# The inserted row picks up the header row's bold/bordered style for its
# interior cells -- clear that back to the plain data-row look first.
$ws.Range("B2:D2").ClearFormats()

# Column A carries the "s=2" index style used throughout this sheet; copy
# it explicitly since the blank row inserted above row 3 didn't inherit it.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 2.31

# Renumber the index column for the rows that shifted down.
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
